$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.256.79"
$ws.Range("E2").Value = "  +1.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.778.19"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.60%  "

# Row 5 - BNB
$ws.Range("D5").Value = "621.59"
$ws.Range("E5").Value = "  +3.84%  "

# Row 6 - Solana
$ws.Range("D6").Value = "165.23"
$ws.Range("E6").Value = "  +1.53%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.776.43"
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.21%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +1.38%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.60%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.43%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  +1.35%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "35.57"
$ws.Range("E14").Value = "  +0.94%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.410.40"
$ws.Range("E15").Value = "  +0.01%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.783.17"
$ws.Range("E16").Value = "  -0.10%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.256.37"
$ws.Range("E17").Value = "  +1.97%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "17.65"
$ws.Range("E18").Value = "  -3.36%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "7.10"
$ws.Range("E19").Value = "  +1.57%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.97%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "468.10"
$ws.Range("E21").Value = "  +2.49%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "9.62"
$ws.Range("E22").Value = "  +0.70%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +1.36%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  +5.78%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "83.30"
$ws.Range("E25").Value = "  +0.66%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "12.03"
$ws.Range("E26").Value = "  +1.12%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  +4.17%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "10.05"
$ws.Range("E28").Value = "  +1.78%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.00%  "

# Row 30 - WrappedeETH
$ws.Range("D30").Value = "3.925.95"
$ws.Range("E30").Value = "  -0.06%  "

# Row 31 - was PancakeSwap, now ImmutableX
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  +1.79%  "

# Row 32 - was ImmutableX, now PancakeSwap
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.66"
$ws.Range("E32").Value = "  +3.78%  "

# Row 33 - NEARProtocol
$ws.Range("D33").Value = "7.31"
$ws.Range("E33").Value = "  +1.28%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "28.79"
$ws.Range("E34").Value = "  -0.39%  "

# Row 35 - Binance-PegBSC-USD
$ws.Range("E35").Value = "  -0.04%  "

# Row 36 - RenzoRestakedETH
$ws.Range("D36").Value = "3.728.54"
$ws.Range("E36").Value = "  +0.09%  "

# Row 37 - Aptos
$ws.Range("D37").Value = "8.98"
$ws.Range("E37").Value = "  +1.03%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +14.29%  "

# Row 39 - Hedera
$ws.Range("E39").Value = "  +3.23%  "

# Row 40 - dogwifhat
$ws.Range("D40").Value = "3.41"
$ws.Range("E40").Value = "  +8.72%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +0.80%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "0.967"
$ws.Range("E42").Value = "  -1.31%  "

# Row 43 - FirstDigitalUSD
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  +0.01%  "

# Row 45 - TheGraph
$ws.Range("D45").Value = "0.299"

# Row 46 - was Monero, now Arweave
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "43.39"
$ws.Range("E46").Value = "  +0.55%  "

# Row 47 - was Arweave, now Monero
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "153.15"
$ws.Range("E47").Value = "  +0.63%  "

# Row 48 - OKB
$ws.Range("D48").Value = "46.77"

# Row 49 - Stacks
$ws.Range("E49").Value = "  +3.94%  "

# Row 50 - Cosmos
$ws.Range("D50").Value = "8.42"
$ws.Range("E50").Value = "  +1.79%  "

# Row 51 - ONDO
$ws.Range("E51").Value = "  +0.60%  "
